# Haver pull Q3 revision
# - Revises previously-published monthly values for rows 132-151, 186-210, and 619-622
#   (Haver re-pull picked up small benchmark revisions in those months).
# - Appends a new monthly observation row (row 623, date serial 44500 = 2021-10-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised values: rows 132-151 ---
$ws.Cells.Item(132, 8).Value = 935099
$ws.Cells.Item(132, 13).Value = 14066646
$ws.Cells.Item(133, 13).Value = 14469098
$ws.Cells.Item(134, 13).Value = 14604733
$ws.Cells.Item(135, 13).Value = 14690170
$ws.Cells.Item(136, 13).Value = 14866629
$ws.Cells.Item(137, 13).Value = 14861148
$ws.Cells.Item(138, 13).Value = 14670539
$ws.Cells.Item(139, 13).Value = 14466764
$ws.Cells.Item(140, 4).Value = 823662
$ws.Cells.Item(140, 13).Value = 14131544
$ws.Cells.Item(140, 14).Value = 8935498
$ws.Cells.Item(141, 13).Value = 13893841
$ws.Cells.Item(141, 14).Value = 8782013
$ws.Cells.Item(142, 13).Value = 13749747
$ws.Cells.Item(142, 14).Value = 8776281
$ws.Cells.Item(143, 13).Value = 13621616
$ws.Cells.Item(143, 14).Value = 8836723
$ws.Cells.Item(144, 14).Value = 9042531
$ws.Cells.Item(145, 14).Value = 9395151
$ws.Cells.Item(146, 14).Value = 9716015
$ws.Cells.Item(147, 14).Value = 9916846
$ws.Cells.Item(148, 14).Value = 10176479
$ws.Cells.Item(149, 14).Value = 10402060
$ws.Cells.Item(150, 14).Value = 10586460
$ws.Cells.Item(151, 14).Value = 10834373

# --- Revised values: rows 186-210 ---
$ws.Cells.Item(186, 4).Value = 590411
$ws.Cells.Item(186, 8).Value = 1209218
$ws.Cells.Item(186, 13).Value = 13890008
$ws.Cells.Item(186, 14).Value = 8211216
$ws.Cells.Item(187, 13).Value = 13966241
$ws.Cells.Item(187, 14).Value = 8231207
$ws.Cells.Item(188, 13).Value = 14177474
$ws.Cells.Item(188, 14).Value = 8310993
$ws.Cells.Item(189, 13).Value = 14260782
$ws.Cells.Item(189, 14).Value = 8315353
$ws.Cells.Item(190, 13).Value = 14423109
$ws.Cells.Item(190, 14).Value = 8370934
$ws.Cells.Item(191, 4).Value = 573527
$ws.Cells.Item(191, 8).Value = 1072983
$ws.Cells.Item(191, 13).Value = 14528189
$ws.Cells.Item(191, 14).Value = 8368435
$ws.Cells.Item(192, 13).Value = 14532964
$ws.Cells.Item(192, 14).Value = 8290009
$ws.Cells.Item(193, 13).Value = 14761487
$ws.Cells.Item(193, 14).Value = 8365101
$ws.Cells.Item(194, 13).Value = 14844168
$ws.Cells.Item(194, 14).Value = 8387988
$ws.Cells.Item(195, 13).Value = 14889099
$ws.Cells.Item(195, 14).Value = 8325460
$ws.Cells.Item(196, 13).Value = 14967780
$ws.Cells.Item(196, 14).Value = 8319956
$ws.Cells.Item(197, 13).Value = 15089925
$ws.Cells.Item(197, 14).Value = 8341420
$ws.Cells.Item(198, 13).Value = 15141431
$ws.Cells.Item(198, 14).Value = 8330287
$ws.Cells.Item(199, 5).Value = 214625
$ws.Cells.Item(199, 13).Value = 15308975
$ws.Cells.Item(199, 14).Value = 8356535
$ws.Cells.Item(199, 15).Value = 2587930
$ws.Cells.Item(200, 13).Value = 15440585
$ws.Cells.Item(200, 14).Value = 8350554
$ws.Cells.Item(200, 15).Value = 2599983
$ws.Cells.Item(201, 13).Value = 15510359
$ws.Cells.Item(201, 14).Value = 8368297
$ws.Cells.Item(201, 15).Value = 2605562
$ws.Cells.Item(202, 13).Value = 15673669
$ws.Cells.Item(202, 14).Value = 8381854
$ws.Cells.Item(202, 15).Value = 2635474
$ws.Cells.Item(203, 15).Value = 2650971
$ws.Cells.Item(204, 15).Value = 2658793
$ws.Cells.Item(205, 15).Value = 2687342
$ws.Cells.Item(206, 15).Value = 2682396
$ws.Cells.Item(207, 15).Value = 2694747
$ws.Cells.Item(208, 15).Value = 2719330
$ws.Cells.Item(209, 15).Value = 2707772
$ws.Cells.Item(210, 15).Value = 2680337

# --- Revised values: rows 619-622 ---
$ws.Cells.Item(619, 2).Value = 2301654
$ws.Cells.Item(619, 4).Value = 658602
$ws.Cells.Item(619, 5).Value = 334050
$ws.Cells.Item(619, 6).Value = 13917720
$ws.Cells.Item(619, 7).Value = 12419612
$ws.Cells.Item(619, 8).Value = 3987205
$ws.Cells.Item(619, 9).Value = 344.61
$ws.Cells.Item(619, 11).Value = 27.55
$ws.Cells.Item(619, 12).Value = 321.48
$ws.Cells.Item(619, 13).Value = 97722359
$ws.Cells.Item(619, 14).Value = 11763945
$ws.Cells.Item(619, 15).Value = 13081476
$ws.Cells.Item(620, 2).Value = 2049168
$ws.Cells.Item(620, 4).Value = 590771
$ws.Cells.Item(620, 5).Value = 286424
$ws.Cells.Item(620, 6).Value = 12847808
$ws.Cells.Item(620, 7).Value = 10786046
$ws.Cells.Item(620, 8).Value = 3457880
$ws.Cells.Item(620, 10).Value = 40.05
$ws.Cells.Item(620, 11).Value = 26.51
$ws.Cells.Item(620, 12).Value = 326.54
$ws.Cells.Item(620, 13).Value = 82587126
$ws.Cells.Item(620, 14).Value = 10198603
$ws.Cells.Item(620, 15).Value = 12418076
$ws.Cells.Item(621, 2).Value = 1504461
$ws.Cells.Item(621, 4).Value = 459840
$ws.Cells.Item(621, 5).Value = 274281
$ws.Cells.Item(621, 6).Value = 12787374
$ws.Cells.Item(621, 7).Value = 10604782
$ws.Cells.Item(621, 8).Value = 3447449
$ws.Cells.Item(621, 11).Value = 24.04
$ws.Cells.Item(621, 12).Value = 333.27
$ws.Cells.Item(621, 13).Value = 69691423
$ws.Cells.Item(621, 14).Value = 9343638
$ws.Cells.Item(621, 15).Value = 11625864
$ws.Cells.Item(622, 2).Value = 1408523
$ws.Cells.Item(622, 4).Value = 480304
$ws.Cells.Item(622, 5).Value = 322658
$ws.Cells.Item(622, 7).Value = 8292609
$ws.Cells.Item(622, 8).Value = 2727366
$ws.Cells.Item(622, 9).Value = 348.85
$ws.Cells.Item(622, 10).Value = 29.13
$ws.Cells.Item(622, 11).Value = 21
$ws.Cells.Item(622, 12).Value = 339.98
$ws.Cells.Item(622, 13).Value = 59033301
$ws.Cells.Item(622, 14).Value = 8912235
$ws.Cells.Item(622, 15).Value = 8865458

# --- New row 623: 2021-10 monthly observation ---
$ws.Cells.Item(623, 1).Value = 44500
$ws.Cells.Item(623, 2).Value = 1238505
$ws.Cells.Item(623, 3).Value = 2945.733333333333
$ws.Cells.Item(623, 4).Value = 360189
$ws.Cells.Item(623, 5).Value = 274338
$ws.Cells.Item(623, 6).Value = 9127049
$ws.Cells.Item(623, 7).Value = 6967141
$ws.Cells.Item(623, 8).Value = 2292105
$ws.Cells.Item(623, 9).Value = 347.35
$ws.Cells.Item(623, 10).Value = 34.8
$ws.Cells.Item(623, 11).Value = 19.49
$ws.Cells.Item(623, 12).Value = 344.19
$ws.Cells.Item(623, 13).Value = 52998131
$ws.Cells.Item(623, 14).Value = 8511303
$ws.Cells.Item(623, 15).Value = 6705815
$ws.Cells.Item(623, 1).NumberFormat = "yyyy-mm-dd"
